$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EMIDATA")

# ----------------------------------------------------------------------
# 1) Insert a new row before row 3 (a split/duplicate of row 2's EMI entry)
#    while the sheet still has its original A:N layout, then copy the
#    formatting of row 2 down onto the new row 3 so every column keeps the
#    same per-column style (bordered plain / currency-format / etc).
# ----------------------------------------------------------------------
$ws.Rows("3:3").Insert()
$ws.Range("A2:N2").Copy()
$ws.Range("A3:N3").PasteSpecial(-4122)

# New row 3 values (old-layout columns): month=24, interestRate=4.4, price=100000
$ws.Range("B3").Value = 24
$ws.Range("C3").Value = 4.4
$ws.Range("D3").Value = 100000

# ----------------------------------------------------------------------
# 2) Insert a new column before column B ("type"), shifting old B:N -> C:O
#    for every row (including the row we just inserted).
# ----------------------------------------------------------------------
$ws.Columns("B:B").Insert()

# ----------------------------------------------------------------------
# 3) New column A = "ID" header + sequential row numbers 1-6
# ----------------------------------------------------------------------
$ws.Range("A1").Value = "ID"
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6

# ----------------------------------------------------------------------
# 4) New column B header = "type" (the old column-A header), with the
#    per-row group id values.
# ----------------------------------------------------------------------
$ws.Range("B1").Value = "type"
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 4
$ws.Range("B7").Value = 5

Write-Host "core edits done"
